# Update dolar_bevsa_uyu historical data: prepend 3 new trading days
# (2026-02-11, 2026-02-12, 2026-02-13) and drop the oldest day (2026-01-12)
# to keep the rolling window, shifting all other rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above the current first data row (row 3),
# pushing existing data rows down.
$ws.Range("A3:H5").EntireRow.Insert()

# New rows (most recent first) to place at the top of the data block.
$newRows = @(
    @(46066, 38.841, 38.841, 38.87, 38.83, 38.95, 42, 21000000),
    @(46065, 38.731, 38.731, 38.75, 38.74, 38.76, 92, 46600000),
    @(46064, 38.537, 38.537, 38.65, 38.62, 38.75, 46, 29000000)
)

# Column A uses a custom date/time number format; grab it from the first
# untouched original data row (now pushed down to row 6) before the loop,
# since rows 3-5 themselves start out with the default "General" format.
$dateFormat = $ws.Cells.Item(6, 1).NumberFormat

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 3 + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 1).NumberFormat = $dateFormat
}

# After the insert, the old last data row (previously row 24, the oldest
# entry 2026-01-12) now sits at row 27, beyond the rolling window we keep.
# Remove it so the table only spans rows 1:26.
$ws.Range("A27:H27").EntireRow.Delete()
